$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Volume/date header text edits (rich-text runs, edited in place via Characters to preserve per-run formatting) ---

# A8: "Volume 32   Number  24" -> "...  25" (last run "24" -> "25")
$volCell = $ws.Range("A8")
$volText = $volCell.Value()
$numIdx = $volText.LastIndexOf("24")
$volChars = $volCell.Characters($numIdx + 1, 2)
$volChars.Text = "25"

# C9: "Report Covering the Week  6/9/2025  Through  6/15/2025"
#     -> "...  6/16/2025  Through  6/22/2025"
$dateCell = $ws.Range("C9")
$dateText = $dateCell.Value()
$i1 = $dateText.IndexOf("6/9/2025")
$i2 = $dateText.IndexOf("6/15/2025")
# Replace the later occurrence first so the earlier index stays valid
$c2 = $dateCell.Characters($i2 + 1, 9)
$c2.Text = "6/22/2025"
$c1 = $dateCell.Characters($i1 + 1, 8)
$c1.Text = "6/16/2025"

# --- Step 1: fix cells whose type changes (number<->text) using donor-cell Copy to preserve exact style ---
$ws.Range("C15").Copy($ws.Range("D15"))  # -> text "0" (style 13, shared string 20)
$ws.Range("M25").Copy($ws.Range("E15"))  # -> text "***.*" (style 13, shared string 21)
$ws.Range("C15").Copy($ws.Range("D27"))  # -> text "0" (style 13, shared string 20)
$ws.Range("M25").Copy($ws.Range("E27"))  # -> text "***.*" (style 13, shared string 21)
$ws.Range("C15").Copy($ws.Range("G29"))  # -> text "0" (style 13, shared string 20)
$ws.Range("M25").Copy($ws.Range("H29"))  # -> text "***.*" (style 13, shared string 21)
$ws.Range("C15").Copy($ws.Range("G30"))  # -> text "0" (style 13, shared string 20)
$ws.Range("M25").Copy($ws.Range("H30"))  # -> text "***.*" (style 13, shared string 21)
$ws.Range("C15").Copy($ws.Range("D31"))  # -> text "0" (style 13, shared string 20)
$ws.Range("M25").Copy($ws.Range("E31"))  # -> text "***.*" (style 13, shared string 21)
$ws.Range("C22").Copy($ws.Range("C23"))  # -> numeric (style 14)
$ws.Range("C23").Value = 2

# --- Step 2: plain numeric value updates ---
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 350
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 103
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = -14.166666666666
$ws.Range("L16").Value = 60.9375
$ws.Range("M16").Value = 6.185567010309
$ws.Range("N16").Value = -81.441441441441
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -81.25
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 13.793103448275
$ws.Range("I17").Value = 166
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 5.732484076433
$ws.Range("L17").Value = 44.347826086956
$ws.Range("M17").Value = 191.228070175439
$ws.Range("N17").Value = -12.631578947368
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -26.086956521739
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = 0.9900990099
$ws.Range("L18").Value = -8.928571428571
$ws.Range("M18").Value = 92.452830188679
$ws.Range("N18").Value = -71.900826446281
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -10.526315789473
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -34.722222222222
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = -7.758620689655
$ws.Range("L19").Value = -1.834862385321
$ws.Range("M19").Value = 56.585365853658
$ws.Range("N19").Value = -32.278481012658
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 43
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = 34.375
$ws.Range("L20").Value = 26.470588235294
$ws.Range("M20").Value = 104.761904761905
$ws.Range("N20").Value = -87.125748502994
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 135
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = -14.012738853503
$ws.Range("I21").Value = 743
$ws.Range("J21").Value = 763
$ws.Range("K21").Value = -2.62123197903
$ws.Range("L21").Value = 13.435114503816
$ws.Range("M21").Value = 70.804597701149
$ws.Range("N21").Value = -61.462655601659
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 3.703703703703
$ws.Range("M22").Value = -15.151515151515
$ws.Range("C23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = -22.727272727272
$ws.Range("L23").Value = -22.727272727272
$ws.Range("M23").Value = 88.888888888888
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 171
$ws.Range("H24").Value = -33.91812865497
$ws.Range("I24").Value = 770
$ws.Range("J24").Value = 1015
$ws.Range("K24").Value = -24.137931034482
$ws.Range("L24").Value = -19.032597266035
$ws.Range("M24").Value = 14.754098360655
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = -19.354838709677
$ws.Range("F25").Value = 99
$ws.Range("G25").Value = 151
$ws.Range("H25").Value = -34.437086092715
$ws.Range("I25").Value = 588
$ws.Range("J25").Value = 918
$ws.Range("K25").Value = -35.9477124183
$ws.Range("L25").Value = -30.904817861339
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = -27.027027027027
$ws.Range("I26").Value = 195
$ws.Range("J26").Value = 230
$ws.Range("K26").Value = -15.217391304347
$ws.Range("L26").Value = 12.068965517241
$ws.Range("M26").Value = 10.795454545454
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 42
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = 75
$ws.Range("L28").Value = 82.608695652173
$ws.Range("L29").Value = -83.333333333333
$ws.Range("M29").Value = -75
$ws.Range("L30").Value = -80
$ws.Range("M30").Value = -50
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 6
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = -45.454545454545
